# Temporary fix for Broly
# Updates the "Summon Rating" (column D) values for the affected rows.
# Values are written with a leading apostrophe so Excel stores them as text
# (matching the inlineStr cells in the workbook) rather than coercing to a
# number, then the style is reset to "Normal" so the quote-prefix formatting
# introduced by the apostrophe does not linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'6.37632562457711"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'14.637989362345085"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'15.810561449629022"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'10.269321654999558"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'2.427615066640635"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'2.546548380998475"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'4.2015699487283396"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'2.34820770565918"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'5.473564005455468"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Value = "'6.430840546867309"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'5.074897300499384"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'2.1098563988051193"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'12.68429271447673"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.7712496821925272"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'6.8115696071666"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'10.851483964559458"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'5.68214487220917"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "'1.7281619174213612"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'7.050245004990942"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'8.530209824755405"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'14.06608577073376"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'7.584452299997466"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'6.349527232251154"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'1.982134295944779"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'7.581071574418942"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'4.386088009954863"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'5.342219454354002"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'3.9315233387901625"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'18.7355189998293"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Value = "'12.72436054811363"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'14.077351771980636"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'4.548024109543142"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'4.968394052950249"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'1.00771638044802"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Value = "'7.018863192545318"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Value = "'14.30659213177179"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'2.0249631229255742"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'5.32640328305607"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").Value = "'5.40522224533119"
$ws.Range("D51").Style = "Normal"
$ws.Range("D52").Value = "'3.737871855704123"
$ws.Range("D52").Style = "Normal"
$ws.Range("D53").Value = "'5.918653708507097"
$ws.Range("D53").Style = "Normal"
$ws.Range("D54").Value = "'16.84424457928631"
$ws.Range("D54").Style = "Normal"
$ws.Range("D56").Value = "'2.746403364385035"
$ws.Range("D56").Style = "Normal"
$ws.Range("D58").Value = "'3.44498357515081"
$ws.Range("D58").Style = "Normal"
$ws.Range("D59").Value = "'2.544789606884329"
$ws.Range("D59").Style = "Normal"
$ws.Range("D62").Value = "'2.1355322402961505"
$ws.Range("D62").Style = "Normal"
$ws.Range("D63").Value = "'0.7349282260217898"
$ws.Range("D63").Style = "Normal"
$ws.Range("D64").Value = "'14.79052575228991"
$ws.Range("D64").Style = "Normal"
$ws.Range("D65").Value = "'20.202554779682536"
$ws.Range("D65").Style = "Normal"
$ws.Range("D66").Value = "'4.137871572923594"
$ws.Range("D66").Style = "Normal"
$ws.Range("D67").Value = "'75.53734087275534"
$ws.Range("D67").Style = "Normal"
$ws.Range("D68").Value = "'2.0550087424665264"
$ws.Range("D68").Style = "Normal"
$ws.Range("D69").Value = "'4.6904200117001515"
$ws.Range("D69").Style = "Normal"
$ws.Range("D70").Value = "'6.185120714599598"
$ws.Range("D70").Style = "Normal"
$ws.Range("D71").Value = "'4.712910075774253"
$ws.Range("D71").Style = "Normal"
$ws.Range("D72").Value = "'5.98917452111926"
$ws.Range("D72").Style = "Normal"
$ws.Range("D73").Value = "'9.91049401614133"
$ws.Range("D73").Style = "Normal"
$ws.Range("D74").Value = "'6.179908529279587"
$ws.Range("D74").Style = "Normal"
$ws.Range("D75").Value = "'4.231544968281098"
$ws.Range("D75").Style = "Normal"
$ws.Range("D76").Value = "'7.230325301975195"
$ws.Range("D76").Style = "Normal"
$ws.Range("D77").Value = "'10.335636177509354"
$ws.Range("D77").Style = "Normal"
$ws.Range("D78").Value = "'5.753596095761784"
$ws.Range("D78").Style = "Normal"
$ws.Range("D79").Value = "'0.7057354171204595"
$ws.Range("D79").Style = "Normal"
$ws.Range("D80").Value = "'22.95388476217289"
$ws.Range("D80").Style = "Normal"
$ws.Range("D81").Value = "'6.6857474162912816"
$ws.Range("D81").Style = "Normal"
$ws.Range("D82").Value = "'96.78167349910503"
$ws.Range("D82").Style = "Normal"
$ws.Range("D83").Value = "'9.834949930463743"
$ws.Range("D83").Style = "Normal"
$ws.Range("D84").Value = "'43.48660725790063"
$ws.Range("D84").Style = "Normal"
$ws.Range("D85").Value = "'12.45242338739863"
$ws.Range("D85").Style = "Normal"
$ws.Range("D86").Value = "'5.792870089277115"
$ws.Range("D86").Style = "Normal"
$ws.Range("D87").Value = "'2.1324162546209573"
$ws.Range("D87").Style = "Normal"
$ws.Range("D88").Value = "'1.759583149849845"
$ws.Range("D88").Style = "Normal"
$ws.Range("D89").Value = "'2.4108123467650544"
$ws.Range("D89").Style = "Normal"
$ws.Range("D90").Value = "'1.308197811798181"
$ws.Range("D90").Style = "Normal"
$ws.Range("D94").Value = "'9.747642269411108"
$ws.Range("D94").Style = "Normal"
$ws.Range("D97").Value = "'2.7052132254043"
$ws.Range("D97").Style = "Normal"
$ws.Range("D98").Value = "'24.855748049467188"
$ws.Range("D98").Style = "Normal"
$ws.Range("D100").Value = "'5.869178962246176"
$ws.Range("D100").Style = "Normal"
$ws.Range("D101").Value = "'40.8951168320198"
$ws.Range("D101").Style = "Normal"
$ws.Range("D102").Value = "'7.259490745487004"
$ws.Range("D102").Style = "Normal"
$ws.Range("D103").Value = "'22.1004788736988"
$ws.Range("D103").Style = "Normal"
$ws.Range("D104").Value = "'1.3455717716039233"
$ws.Range("D104").Style = "Normal"
$ws.Range("D106").Value = "'7.114400308121961"
$ws.Range("D106").Style = "Normal"
$ws.Range("D107").Value = "'2.64772765297398"
$ws.Range("D107").Style = "Normal"
$ws.Range("D110").Value = "'5.764061235966218"
$ws.Range("D110").Style = "Normal"
$ws.Range("D111").Value = "'21.21411771073802"
$ws.Range("D111").Style = "Normal"
$ws.Range("D114").Value = "'1.1704637820302892"
$ws.Range("D114").Style = "Normal"
$ws.Range("D117").Value = "'0.43758979377584517"
$ws.Range("D117").Style = "Normal"
$ws.Range("D118").Value = "'2.3119148846580018"
$ws.Range("D118").Style = "Normal"
$ws.Range("D119").Value = "'5.071622736239579"
$ws.Range("D119").Style = "Normal"
$ws.Range("D120").Value = "'9.164076992184583"
$ws.Range("D120").Style = "Normal"
$ws.Range("D121").Value = "'2.0515469862802007"
$ws.Range("D121").Style = "Normal"
$ws.Range("D122").Value = "'3.770307462668546"
$ws.Range("D122").Style = "Normal"
$ws.Range("D124").Value = "'0.37597752388730604"
$ws.Range("D124").Style = "Normal"
$ws.Range("D125").Value = "'4.613173510921976"
$ws.Range("D125").Style = "Normal"
$ws.Range("D128").Value = "'0.2938871194131131"
$ws.Range("D128").Style = "Normal"
$ws.Range("D135").Value = "'2.7790200213606404"
$ws.Range("D135").Style = "Normal"
$ws.Range("D136").Value = "'1.5429381724920674"
$ws.Range("D136").Style = "Normal"
$ws.Range("D139").Value = "'1.5371379841690436"
$ws.Range("D139").Style = "Normal"
$ws.Range("D140").Value = "'0.5766293490807008"
$ws.Range("D140").Style = "Normal"
$ws.Range("D141").Value = "'1.866942179953085"
$ws.Range("D141").Style = "Normal"
$ws.Range("D145").Value = "'1.8869496557045216"
$ws.Range("D145").Style = "Normal"
$ws.Range("D147").Value = "'2.5967256189951096"
$ws.Range("D147").Style = "Normal"
$ws.Range("D150").Value = "'6.061535597579528"
$ws.Range("D150").Style = "Normal"
$ws.Range("D151").Value = "'3.7598169538776336"
$ws.Range("D151").Style = "Normal"
$ws.Range("D152").Value = "'3.977102310729566"
$ws.Range("D152").Style = "Normal"
$ws.Range("D154").Value = "'2.364702110842875"
$ws.Range("D154").Style = "Normal"
$ws.Range("D155").Value = "'0.8345435242855277"
$ws.Range("D155").Style = "Normal"
$ws.Range("D157").Value = "'2.5903710691932207"
$ws.Range("D157").Style = "Normal"
$ws.Range("D158").Value = "'16.61159760445597"
$ws.Range("D158").Style = "Normal"
$ws.Range("D159").Value = "'9.167574413979871"
$ws.Range("D159").Style = "Normal"
$ws.Range("D160").Value = "'19.79604225553276"
$ws.Range("D160").Style = "Normal"
$ws.Range("D161").Value = "'9.60809745399905"
$ws.Range("D161").Style = "Normal"
$ws.Range("D162").Value = "'22.372326922838695"
$ws.Range("D162").Style = "Normal"
$ws.Range("D163").Value = "'13.42000546627794"
$ws.Range("D163").Style = "Normal"
$ws.Range("D164").Value = "'2.9633351177829974"
$ws.Range("D164").Style = "Normal"
$ws.Range("D166").Value = "'17.31134230584965"
$ws.Range("D166").Style = "Normal"
$ws.Range("D171").Value = "'6.383774019915773"
$ws.Range("D171").Style = "Normal"
$ws.Range("D173").Value = "'4.40903942137464"
$ws.Range("D173").Style = "Normal"
$ws.Range("D174").Value = "'1.7840390227554954"
$ws.Range("D174").Style = "Normal"
$ws.Range("D185").Value = "'2.3606454353698547"
$ws.Range("D185").Style = "Normal"
$ws.Range("D186").Value = "'0.696097880043544"
$ws.Range("D186").Style = "Normal"
$ws.Range("D187").Value = "'5.811252869765881"
$ws.Range("D187").Style = "Normal"
$ws.Range("D190").Value = "'2.4351914080931376"
$ws.Range("D190").Style = "Normal"
$ws.Range("D191").Value = "'2.198737239288866"
$ws.Range("D191").Style = "Normal"
$ws.Range("D193").Value = "'0.28384418931991534"
$ws.Range("D193").Style = "Normal"
$ws.Range("D194").Value = "'0.36272899202234754"
$ws.Range("D194").Style = "Normal"
$ws.Range("D199").Value = "'0.40680695568551983"
$ws.Range("D199").Style = "Normal"
$ws.Range("D204").Value = "'2.0631122721488673"
$ws.Range("D204").Style = "Normal"
$ws.Range("D205").Value = "'1.3778231213128702"
$ws.Range("D205").Style = "Normal"
$ws.Range("D206").Value = "'1.1513359738702909"
$ws.Range("D206").Style = "Normal"
$ws.Range("D207").Value = "'0.35994806169069243"
$ws.Range("D207").Style = "Normal"
$ws.Range("D210").Value = "'0.8570983944895296"
$ws.Range("D210").Style = "Normal"
$ws.Range("D212").Value = "'12.489790205676222"
$ws.Range("D212").Style = "Normal"
$ws.Range("D215").Value = "'82.33539619162269"
$ws.Range("D215").Style = "Normal"
$ws.Range("D217").Value = "'17.8598189530847"
$ws.Range("D217").Style = "Normal"
$ws.Range("D218").Value = "'17.31511670490435"
$ws.Range("D218").Style = "Normal"
$ws.Range("D219").Value = "'95.77899130394034"
$ws.Range("D219").Style = "Normal"
$ws.Range("D223").Value = "'15.986454877423068"
$ws.Range("D223").Style = "Normal"
